$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (32 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3335.2942
$ws.Range("I106").Value = 2914.2856
$ws.Range("J106").Value = 5300
$ws.Range("K106").Value = 2914.2856
$ws.Range("L106").Value = 5300
$ws.Range("M106").Value = -2283.2856
$ws.Range("N106").Value = -6562
$ws.Range("H112").Value = 6265.4517
$ws.Range("I112").Value = 668.75
$ws.Range("J112").Value = 8212.130999999999
$ws.Range("K112").Value = 2006.25
$ws.Range("L112").Value = 24636.393
$ws.Range("M112").Value = -898.25
$ws.Range("N112").Value = -26852.393
$ws.Range("H113").Value = 2590.3333
$ws.Range("I113").Value = 2260.5
$ws.Range("J113").Value = 3250
$ws.Range("K113").Value = 2260.5
$ws.Range("L113").Value = 3250
$ws.Range("M113").Value = 993.5
$ws.Range("N113").Value = -9758
$ws.Range("H132").Value = 1558.5143
$ws.Range("I132").Value = 1038.1
$ws.Range("J132").Value = 4681
$ws.Range("K132").Value = 3114.3
$ws.Range("L132").Value = 14043
$ws.Range("M132").Value = -584.2999999999997
$ws.Range("N132").Value = -19103
$ws.Range("H137").Value = 793.6383
$ws.Range("I137").Value = 657.4375
$ws.Range("K137").Value = 1972.3125
$ws.Range("M137").Value = 577.6875

# ---- Sheet: ARM (39 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1446.1333
$ws.Range("I61").Value = 1322.4615
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 1322.4615
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -1110.4615
$ws.Range("N61").Value = -2674
$ws.Range("H74").Value = 1264.138
$ws.Range("I74").Value = 1387.9166
$ws.Range("J74").Value = 670
$ws.Range("K74").Value = 1387.9166
$ws.Range("L74").Value = 670
$ws.Range("M74").Value = -513.9166
$ws.Range("N74").Value = -2418
$ws.Range("H77").Value = 1264.138
$ws.Range("I77").Value = 1387.9166
$ws.Range("J77").Value = 670
$ws.Range("K77").Value = 6939.583000000001
$ws.Range("L77").Value = 3350
$ws.Range("M77").Value = -2571.583000000001
$ws.Range("N77").Value = -12086
$ws.Range("H102").Value = 4800
$ws.Range("I102").Value = 3866.6667
$ws.Range("K102").Value = 3866.6667
$ws.Range("M102").Value = -2244.6667
$ws.Range("H132").Value = 1773.2084
$ws.Range("I132").Value = 1325.1904
$ws.Range("J132").Value = 4909.3335
$ws.Range("K132").Value = 3975.5712
$ws.Range("L132").Value = 14728.0005
$ws.Range("M132").Value = -1445.5712
$ws.Range("N132").Value = -19788.0005
$ws.Range("H136").Value = 1446.1333
$ws.Range("I136").Value = 1322.4615
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 3967.3845
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -1417.3845
$ws.Range("N136").Value = -11850

# ---- Sheet: BSM (11 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1713.2222
$ws.Range("I107").Value = 1316.375
$ws.Range("J107").Value = 4888
$ws.Range("K107").Value = 1316.375
$ws.Range("L107").Value = 4888
$ws.Range("M107").Value = 603.625
$ws.Range("N107").Value = -8728
$ws.Range("H134").Value = 1349.5264
$ws.Range("I134").Value = 1005.81036
$ws.Range("K134").Value = 3017.43108
$ws.Range("M134").Value = -482.4310799999998

# ---- Sheet: CRP (63 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2981.6394
$ws.Range("I31").Value = 1546.2693
$ws.Range("J31").Value = 11274.889
$ws.Range("K31").Value = 1546.2693
$ws.Range("L31").Value = 11274.889
$ws.Range("M31").Value = -1251.2693
$ws.Range("N31").Value = -11864.889
$ws.Range("H34").Value = 2981.6394
$ws.Range("I34").Value = 1546.2693
$ws.Range("J34").Value = 11274.889
$ws.Range("K34").Value = 1546.2693
$ws.Range("L34").Value = 11274.889
$ws.Range("M34").Value = -1344.2693
$ws.Range("N34").Value = -11678.889
$ws.Range("H58").Value = 1649.0741
$ws.Range("I58").Value = 1456.5625
$ws.Range("J58").Value = 1929.091
$ws.Range("K58").Value = 1456.5625
$ws.Range("L58").Value = 1929.091
$ws.Range("M58").Value = -1253.5625
$ws.Range("N58").Value = -2335.091
$ws.Range("H99").Value = 1361.1842
$ws.Range("I99").Value = 1268.7
$ws.Range("J99").Value = 1708
$ws.Range("K99").Value = 1268.7
$ws.Range("L99").Value = 1708
$ws.Range("M99").Value = 229.3
$ws.Range("N99").Value = -4704
$ws.Range("H105").Value = 1622
$ws.Range("I105").Value = 1231.25
$ws.Range("J105").Value = 2068.5715
$ws.Range("K105").Value = 1231.25
$ws.Range("L105").Value = 2068.5715
$ws.Range("M105").Value = 515.75
$ws.Range("N105").Value = -5562.5715
$ws.Range("H126").Value = 1361.1842
$ws.Range("I126").Value = 1268.7
$ws.Range("J126").Value = 1708
$ws.Range("K126").Value = 3806.1
$ws.Range("L126").Value = 5124
$ws.Range("M126").Value = -1336.1
$ws.Range("N126").Value = -10064
$ws.Range("H132").Value = 1992.8667
$ws.Range("I132").Value = 1665.3
$ws.Range("J132").Value = 2648
$ws.Range("K132").Value = 4995.9
$ws.Range("L132").Value = 7944
$ws.Range("M132").Value = -2465.9
$ws.Range("N132").Value = -13004
$ws.Range("H134").Value = 1713.5625
$ws.Range("I134").Value = 1864.1666
$ws.Range("J134").Value = 1261.75
$ws.Range("K134").Value = 5592.4998
$ws.Range("L134").Value = 3785.25
$ws.Range("M134").Value = -3057.4998
$ws.Range("N134").Value = -8855.25
$ws.Range("H136").Value = 1649.0741
$ws.Range("I136").Value = 1456.5625
$ws.Range("J136").Value = 1929.091
$ws.Range("K136").Value = 4369.6875
$ws.Range("L136").Value = 5787.272999999999
$ws.Range("M136").Value = -1819.6875
$ws.Range("N136").Value = -10887.273

# ---- Sheet: CUL (74 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 19609108
$ws.Range("I129").Value = 676.6667
$ws.Range("J129").Value = 23810914
$ws.Range("K129").Value = 2030.0001
$ws.Range("L129").Value = 71432742
$ws.Range("M129").Value = 2969.9999
$ws.Range("N129").Value = -71442742
$ws.Range("H130").Value = 1515
$ws.Range("I130").Value = 1515
$ws.Range("K130").Value = 4545
$ws.Range("M130").Value = 475
$ws.Range("H131").Value = 7247201.5
$ws.Range("I131").Value = 477.14285
$ws.Range("J131").Value = 8065380
$ws.Range("K131").Value = 1431.42855
$ws.Range("L131").Value = 24196140
$ws.Range("M131").Value = 3608.57145
$ws.Range("N131").Value = -24206220
$ws.Range("H133").Value = 7497.4185
$ws.Range("I133").Value = 7684.143
$ws.Range("J133").Value = 7461.1113
$ws.Range("K133").Value = 23052.429
$ws.Range("L133").Value = 22383.3339
$ws.Range("M133").Value = -17992.429
$ws.Range("N133").Value = -32503.3339
$ws.Range("H134").Value = 3132.889
$ws.Range("I134").Value = 2113.524
$ws.Range("J134").Value = 4560
$ws.Range("K134").Value = 6340.572
$ws.Range("L134").Value = 13680
$ws.Range("M134").Value = -1270.572
$ws.Range("N134").Value = -23820
$ws.Range("H136").Value = 1735.5294
$ws.Range("I136").Value = 778.25
$ws.Range("J136").Value = 4033
$ws.Range("K136").Value = 2334.75
$ws.Range("L136").Value = 12099
$ws.Range("M136").Value = 2765.25
$ws.Range("N136").Value = -22299
$ws.Range("H137").Value = 2138.1143
$ws.Range("I137").Value = 859.2308
$ws.Range("J137").Value = 2893.818
$ws.Range("K137").Value = 2577.6924
$ws.Range("L137").Value = 8681.454000000002
$ws.Range("M137").Value = 2522.3076
$ws.Range("N137").Value = -18881.454
$ws.Range("H138").Value = 1129.591
$ws.Range("I138").Value = 1041.4762
$ws.Range("J138").Value = 2980
$ws.Range("K138").Value = 3124.4286
$ws.Range("L138").Value = 8940
$ws.Range("M138").Value = 2015.5714
$ws.Range("N138").Value = -19220
$ws.Range("H139").Value = 2376.366
$ws.Range("I139").Value = 1088.7333
$ws.Range("J139").Value = 3119.2307
$ws.Range("K139").Value = 3266.199900000001
$ws.Range("L139").Value = 9357.6921
$ws.Range("M139").Value = 1873.800099999999
$ws.Range("N139").Value = -19637.6921
$ws.Range("H140").Value = 7995.8335
$ws.Range("I140").Value = 4005.5557
$ws.Range("J140").Value = 19966.666
$ws.Range("K140").Value = 12016.6671
$ws.Range("L140").Value = 59899.99800000001
$ws.Range("M140").Value = -6836.667099999999
$ws.Range("N140").Value = -70259.99800000001
$ws.Range("H141").Value = 6156.815
$ws.Range("I141").Value = 7718.75
$ws.Range("J141").Value = 5499.1577
$ws.Range("K141").Value = 23156.25
$ws.Range("L141").Value = 16497.4731
$ws.Range("M141").Value = -17976.25
$ws.Range("N141").Value = -26857.4731

# ---- Sheet: GSM (7 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2112.4773
$ws.Range("I132").Value = 1978.2307
$ws.Range("J132").Value = 3159.6
$ws.Range("K132").Value = 5934.6921
$ws.Range("L132").Value = 9478.799999999999
$ws.Range("M132").Value = -3404.6921
$ws.Range("N132").Value = -14538.8

# ---- Sheet: LTW (15 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2478.4138
$ws.Range("J100").Value = 2609.4546
$ws.Range("L100").Value = 2609.4546
$ws.Range("N100").Value = -3691.4546
$ws.Range("H132").Value = 2746.1924
$ws.Range("I132").Value = 2559.0488
$ws.Range("J132").Value = 3443.7273
$ws.Range("K132").Value = 7677.1464
$ws.Range("L132").Value = 10331.1819
$ws.Range("M132").Value = -5147.1464
$ws.Range("N132").Value = -15391.1819
$ws.Range("H136").Value = 2739.8545
$ws.Range("I136").Value = 2129.7
$ws.Range("K136").Value = 6389.099999999999
$ws.Range("M136").Value = -3839.099999999999

# ---- Sheet: WVR (21 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 692.36365
$ws.Range("I107").Value = 702.6
$ws.Range("J107").Value = 590
$ws.Range("K107").Value = 2107.8
$ws.Range("L107").Value = 1770
$ws.Range("M107").Value = -187.8000000000002
$ws.Range("N107").Value = -5610
$ws.Range("H132").Value = 16667679
$ws.Range("I132").Value = 18940274
$ws.Range("J132").Value = 1987
$ws.Range("K132").Value = 56820822
$ws.Range("L132").Value = 5961
$ws.Range("M132").Value = -56818292
$ws.Range("N132").Value = -11021
$ws.Range("H136").Value = 848.8298
$ws.Range("I136").Value = 715.2564
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2145.7692
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = 404.2308000000003
$ws.Range("N136").Value = -9600
